# Natmi following Dr Hou advice
# Update the Bmp2 -> Rgmb LR-pair sheet: add the "ECs" sending cluster so
# every sending cluster (ECs, FAPs, sCs) is paired with every target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Bmp2"
$ws.Range("C2").Value2 = "Rgmb"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 9.163165666666666
$ws.Range("H2").Value2 = 27.489497
$ws.Range("I2").Value2 = 0.2800251397703982
$ws.Range("J2").Value2 = 0.2800251397703982
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 4.367553
$ws.Range("N2").Value2 = 13.102659
$ws.Range("O2").Value2 = 0.05657159077620311
$ws.Range("P2").Value2 = 0.05657159077620311
$ws.Range("Q2").Value2 = 40.020611696947
$ws.Range("R2").Value2 = 360.185505272523
$ws.Range("S2").Value2 = 0.01584146761414004
$ws.Range("T2").Value2 = 0.01584146761414005

# row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Bmp2"
$ws.Range("C3").Value2 = "Rgmb"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 9.163165666666666
$ws.Range("H3").Value2 = 27.489497
$ws.Range("I3").Value2 = 0.2800251397703982
$ws.Range("J3").Value2 = 0.2800251397703982
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 53.45120266666667
$ws.Range("N3").Value2 = 160.353608
$ws.Range("O3").Value2 = 0.6923372340884159
$ws.Range("P3").Value2 = 0.6923372340884159
$ws.Range("Q3").Value2 = 489.7822251172418
$ws.Range("R3").Value2 = 4408.040026055176
$ws.Range("S3").Value2 = 0.1938718307438595
$ws.Range("T3").Value2 = 0.1938718307438596

# row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Bmp2"
$ws.Range("C4").Value2 = "Rgmb"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 9.163165666666666
$ws.Range("H4").Value2 = 27.489497
$ws.Range("I4").Value2 = 0.2800251397703982
$ws.Range("J4").Value2 = 0.2800251397703982
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 19.38524266666667
$ws.Range("N4").Value2 = 58.155728
$ws.Range("O4").Value2 = 0.251091175135381
$ws.Range("P4").Value2 = 0.251091175135381
$ws.Range("Q4").Value2 = 177.6301900432018
$ws.Range("R4").Value2 = 1598.671710388816
$ws.Range("S4").Value2 = 0.0703118414123986
$ws.Range("T4").Value2 = 0.07031184141239862

# row 5
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Bmp2"
$ws.Range("C5").Value2 = "Rgmb"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 12.06704733333333
$ws.Range("H5").Value2 = 36.201142
$ws.Range("I5").Value2 = 0.3687673822623249
$ws.Range("J5").Value2 = 0.3687673822623249
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 4.367553
$ws.Range("N5").Value2 = 13.102659
$ws.Range("O5").Value2 = 0.05657159077620311
$ws.Range("P5").Value2 = 0.05657159077620311
$ws.Range("Q5").Value2 = 52.703468781842
$ws.Range("R5").Value2 = 474.331219036578
$ws.Range("S5").Value2 = 0.0208617574409559
$ws.Range("T5").Value2 = 0.0208617574409559

# row 6
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Bmp2"
$ws.Range("C6").Value2 = "Rgmb"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 12.06704733333333
$ws.Range("H6").Value2 = 36.201142
$ws.Range("I6").Value2 = 0.3687673822623249
$ws.Range("J6").Value2 = 0.3687673822623249
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 53.45120266666667
$ws.Range("N6").Value2 = 160.353608
$ws.Range("O6").Value2 = 0.6923372340884159
$ws.Range("P6").Value2 = 0.6923372340884159
$ws.Range("Q6").Value2 = 644.9981926022596
$ws.Range("R6").Value2 = 5804.983733420336
$ws.Range("S6").Value2 = 0.2553113894575235
$ws.Range("T6").Value2 = 0.2553113894575235

# row 7
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Bmp2"
$ws.Range("C7").Value2 = "Rgmb"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 12.06704733333333
$ws.Range("H7").Value2 = 36.201142
$ws.Range("I7").Value2 = 0.3687673822623249
$ws.Range("J7").Value2 = 0.3687673822623249
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 19.38524266666667
$ws.Range("N7").Value2 = 58.155728
$ws.Range("O7").Value2 = 0.251091175135381
$ws.Range("P7").Value2 = 0.251091175135381
$ws.Range("Q7").Value2 = 233.9226408268196
$ws.Range("R7").Value2 = 2105.303767441376
$ws.Range("S7").Value2 = 0.09259423536384542
$ws.Range("T7").Value2 = 0.09259423536384542

# row 8
$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Bmp2"
$ws.Range("C8").Value2 = "Rgmb"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 11.49244066666667
$ws.Range("H8").Value2 = 34.477322
$ws.Range("I8").Value2 = 0.3512074779672769
$ws.Range("J8").Value2 = 0.351207477967277
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 4.367553
$ws.Range("N8").Value2 = 13.102659
$ws.Range("O8").Value2 = 0.05657159077620311
$ws.Range("P8").Value2 = 0.05657159077620311
$ws.Range("Q8").Value2 = 50.19384371102201
$ws.Range("R8").Value2 = 451.744593399198
$ws.Range("S8").Value2 = 0.01986836572110716
$ws.Range("T8").Value2 = 0.01986836572110716

# row 9
$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Bmp2"
$ws.Range("C9").Value2 = "Rgmb"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 11.49244066666667
$ws.Range("H9").Value2 = 34.477322
$ws.Range("I9").Value2 = 0.3512074779672769
$ws.Range("J9").Value2 = 0.351207477967277
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 53.45120266666667
$ws.Range("N9").Value2 = 160.353608
$ws.Range("O9").Value2 = 0.6923372340884159
$ws.Range("P9").Value2 = 0.6923372340884159
$ws.Range("Q9").Value2 = 614.2847752086418
$ws.Range("R9").Value2 = 5528.562976877776
$ws.Range("S9").Value2 = 0.2431540138870328
$ws.Range("T9").Value2 = 0.2431540138870328

# row 10
$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Bmp2"
$ws.Range("C10").Value2 = "Rgmb"
$ws.Range("D10").Value2 = "sCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 11.49244066666667
$ws.Range("H10").Value2 = 34.477322
$ws.Range("I10").Value2 = 0.3512074779672769
$ws.Range("J10").Value2 = 0.351207477967277
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 19.38524266666667
$ws.Range("N10").Value2 = 58.155728
$ws.Range("O10").Value2 = 0.251091175135381
$ws.Range("P10").Value2 = 0.251091175135381
$ws.Range("Q10").Value2 = 222.7837511556018
$ws.Range("R10").Value2 = 2005.053760400416
$ws.Range("S10").Value2 = 0.088185098359137
$ws.Range("T10").Value2 = 0.08818509835913702
